$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44
$ws.Range("E2").Value = 5

$ws.Range("G4").Select()
